$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.724.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "'3.691.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'672.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "'160.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.499"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("D9").Value = "'0.146"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").Value = "'7.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("D11").Value = "'0.443"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").Value = "'33.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.41%  "
$ws.Range("D14").Value = "'3.700.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "'69.681.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").Value = "'16.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'6.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("D19").Value = "'471.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("D20").Value = "'9.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "'0.648"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "'80.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "'3.840.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("E24").Value = "  +6.99%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'10.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").Value = "'9.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D29").Value = "'1.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  +2.97%  "
$ws.Range("E31").Value = "  +6.03%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "'6.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("D34").Value = "'26.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").Value = "'3.689.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("D36").Value = "'8.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.15%  "
$ws.Range("D37").Value = "'6.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D39").Value = "'2.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.12%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "'0.0908"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").Value = "'175.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.65%  "
$ws.Range("D43").Value = "'0.935"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").Value = "'2.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.88%  "
$ws.Range("D46").Value = "'27.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("D47").Value = "'1.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.84%  "
$ws.Range("D48").Value = "'0.000274"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("D50").Value = "'7.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").Value = "'0.265"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.16%  "
